$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.690.85"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "1.611.32"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.525"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "28.91"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.259"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0609"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0907"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").Value = "1.614.84"
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.566"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.85"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.23%  "
$ws.Range("D16").Value = "29.701.12"
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +12.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("D20").Value = "0.0₃0707"
$ws.Range("E20").Value = "  +2.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.05%  "
$ws.Range("E27").Value = "  +1.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0482"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.14%  "
$ws.Range("E32").Value = "  +0.69%  "
$ws.Range("E33").Value = "  +2.69%  "
$ws.Range("D34").Value = "1.431.63"
$ws.Range("E34").Value = "  +0.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.62"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.36%  "
$ws.Range("E36").Value = "  +1.28%  "
$ws.Range("E37").Value = "  +2.23%  "
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("E39").Value = "  +3.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.557"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0499"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.829"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.41%  "
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "54.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "68.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.997"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +19.97%  "
$ws.Range("E48").Value = "  +3.31%  "
$ws.Range("D49").Value = "1.747.97"
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "87.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.99%  "
$ws.Range("E51").Value = "  -1.10%  "
